$wb = $excel.ActiveWorkbook

# Helper: write a value to a cell. Numeric-looking text (e.g. "7.3", "13.0",
# "11") needs a leading quote so Excel keeps it as text instead of silently
# converting it to a number.
function Set-TextValue {
    param($range, $text)
    if ($text -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# --- CypherOutput sheet: update existing rows 2-5 and add new rows 6-7 ---
$ws = $wb.Worksheets.Item("CypherOutput")

# Row 2 - COTC007B-0412 (Lymphoma / Va / Male)
Set-TextValue $ws.Range("A2") "COTC007B-0412"
Set-TextValue $ws.Range("B2") "COTC007B"
Set-TextValue $ws.Range("C2") "Clinical Trial"
Set-TextValue $ws.Range("D2") "American Staffordshire Terrier"
Set-TextValue $ws.Range("E2") "Lymphoma"
Set-TextValue $ws.Range("F2") "Va"
Set-TextValue $ws.Range("G2") "7.3"
Set-TextValue $ws.Range("H2") "Male"
Set-TextValue $ws.Range("I2") "Yes"

# Row 3 - COTC007B-0301 (Lymphoma / IIIa / Female)
Set-TextValue $ws.Range("A3") "COTC007B-0301"
Set-TextValue $ws.Range("B3") "COTC007B"
Set-TextValue $ws.Range("C3") "Clinical Trial"
Set-TextValue $ws.Range("D3") "American Staffordshire Terrier"
Set-TextValue $ws.Range("E3") "Lymphoma"
Set-TextValue $ws.Range("F3") "IIIa"
Set-TextValue $ws.Range("G3") "5.8"
Set-TextValue $ws.Range("H3") "Female"
Set-TextValue $ws.Range("I3") "Yes"

# Row 4 - COTC007B-0409 (Lymphoma / IIIa / Male)
Set-TextValue $ws.Range("A4") "COTC007B-0409"
Set-TextValue $ws.Range("B4") "COTC007B"
Set-TextValue $ws.Range("C4") "Clinical Trial"
Set-TextValue $ws.Range("D4") "American Staffordshire Terrier"
Set-TextValue $ws.Range("E4") "Lymphoma"
Set-TextValue $ws.Range("F4") "IIIa"
Set-TextValue $ws.Range("G4") "2.4"
Set-TextValue $ws.Range("H4") "Male"
Set-TextValue $ws.Range("I4") "Yes"

# Row 5 - COTC007B-0617 (Lymphoma / III / Female)
Set-TextValue $ws.Range("A5") "COTC007B-0617"
Set-TextValue $ws.Range("B5") "COTC007B"
Set-TextValue $ws.Range("C5") "Clinical Trial"
Set-TextValue $ws.Range("D5") "American Staffordshire Terrier"
Set-TextValue $ws.Range("E5") "Lymphoma"
Set-TextValue $ws.Range("F5") "III"
Set-TextValue $ws.Range("G5") "10.0"
Set-TextValue $ws.Range("H5") "Female"
Set-TextValue $ws.Range("I5") "Yes"

# Row 6 (new) - NCATS-COP01-CCB040254 (Pulmonary Neoplasms / Unknown / Female)
Set-TextValue $ws.Range("A6") "NCATS-COP01-CCB040254"
Set-TextValue $ws.Range("B6") "NCATS-COP01"
Set-TextValue $ws.Range("C6") "Transcriptomics"
Set-TextValue $ws.Range("D6") "American Staffordshire Terrier"
Set-TextValue $ws.Range("E6") "Pulmonary Neoplasms"
Set-TextValue $ws.Range("F6") "Unknown"
Set-TextValue $ws.Range("G6") "13.0"
Set-TextValue $ws.Range("H6") "Female"
Set-TextValue $ws.Range("I6") "Yes"

# Row 7 (new) - GLIOMA01-i_4990 (Glioma / Unknown / Male)
Set-TextValue $ws.Range("A7") "GLIOMA01-i_4990"
Set-TextValue $ws.Range("B7") "GLIOMA01"
Set-TextValue $ws.Range("C7") "Genomics"
Set-TextValue $ws.Range("D7") "American Staffordshire Terrier"
Set-TextValue $ws.Range("E7") "Glioma"
Set-TextValue $ws.Range("F7") "Unknown"
Set-TextValue $ws.Range("G7") "4.0"
Set-TextValue $ws.Range("H7") "Male"
Set-TextValue $ws.Range("I7") "Yes"

# --- StatOutput sheet: update totals row ---
$statWs = $wb.Worksheets.Item("StatOutput")
Set-TextValue $statWs.Range("A2") "11"
Set-TextValue $statWs.Range("B2") "8"
Set-TextValue $statWs.Range("C2") "6"
Set-TextValue $statWs.Range("D2") "3"
